$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values below mirror the scraped coinranking.com snapshot used by the
# GitHub Actions refresh job. D-column price strings are forced to Text via a
# leading apostrophe so values like "1.000" / "308.72" are not auto-converted
# to numbers by Excel (matching the original inlineStr cell type).

$ws.Range("D2").Value = '''26.754.36'
$ws.Range("E2").Value = '  -2.21%  '
$ws.Range("D3").Value = '''1.798.54'
$ws.Range("E3").Value = '  -1.57%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''308.72'
$ws.Range("E5").Value = '  -1.76%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '''0.4591'
$ws.Range("E7").Value = '  +2.51%  '
$ws.Range("D8").Value = '''0.3715'
$ws.Range("E8").Value = '  -1.53%  '
$ws.Range("D9").Value = '''0.07254'
$ws.Range("E9").Value = '  -3.70%  '
$ws.Range("D10").Value = '''0.8561'
$ws.Range("E10").Value = '  -4.14%  '
$ws.Range("D11").Value = '''20.38'
$ws.Range("E11").Value = '  -3.32%  '
$ws.Range("D12").Value = '''1.778.15'
$ws.Range("E12").Value = '  -2.65%  '
$ws.Range("D13").Value = '''5.316'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("D14").Value = '''6.489'
$ws.Range("E14").Value = '  -3.83%  '
$ws.Range("D15").Value = '''0.07035'
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("E16").Value = '  -4.27%  '
$ws.Range("D17").Value = '''1.001'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '''0.000008631'
$ws.Range("E18").Value = '  -2.33%  '
$ws.Range("D19").Value = '''1.0000'
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = '''14.64'
$ws.Range("E20").Value = '  -3.82%  '
$ws.Range("D21").Value = '''26.767.50'
$ws.Range("E21").Value = '  -2.23%  '
$ws.Range("D22").Value = '''5.283'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("D24").Value = '''2.022.41'
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("D25").Value = '''1.904'
$ws.Range("E25").Value = '  -5.00%  '
$ws.Range("D26").Value = '''149.65'
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("D27").Value = '''2.148'
$ws.Range("E27").Value = '  -13.63%  '
$ws.Range("D28").Value = '''18.16'
$ws.Range("E28").Value = '  -2.68%  '
$ws.Range("D29").Value = '''5.213'
$ws.Range("E29").Value = '  -2.96%  '
$ws.Range("D30").Value = '''114.10'
$ws.Range("E30").Value = '  -3.33%  '
$ws.Range("D31").Value = '''0.08864'
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("D32").Value = '''0.7539'
$ws.Range("E32").Value = '  -3.69%  '
$ws.Range("E33").Value = '  -3.83%  '
$ws.Range("D34").Value = '''4.431'
$ws.Range("E34").Value = '  -3.01%  '
$ws.Range("D35").Value = '''2.886'
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = '''0.9995'
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("D38").Value = '''0.01937'
$ws.Range("E38").Value = '  -2.80%  '
$ws.Range("D39").Value = '''0.05212'
$ws.Range("E39").Value = '  -2.36%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.895'
$ws.Range("E40").Value = '  +0.74%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '''2.364'
$ws.Range("E41").Value = '  +2.96%  '
$ws.Range("D42").Value = '''7.174'
$ws.Range("E42").Value = '  -2.91%  '
$ws.Range("D43").Value = '''0.5226'
$ws.Range("E43").Value = '  -2.00%  '
$ws.Range("D44").Value = '''0.1644'
$ws.Range("E44").Value = '  -5.14%  '
$ws.Range("D45").Value = '''8.496'
$ws.Range("E45").Value = '  -3.70%  '
$ws.Range("D46").Value = '''0.4996'
$ws.Range("E46").Value = '  -3.32%  '
$ws.Range("D47").Value = '''10.21'
$ws.Range("E47").Value = '  -5.95%  '
$ws.Range("D48").Value = '''104.08'
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("D49").Value = '''0.9994'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").Value = '''1.644'
$ws.Range("E50").Value = '  -3.97%  '
$ws.Range("D51").Value = '''0.06287'
$ws.Range("E51").Value = '  -1.39%  '
